# The uploaded workbook's header row is being relabeled from the
# descriptive names (id / name / value) to the generic column letters
# (A / B / C). The underlying row data (1/Alice/100, 2/Bob/200,
# 3/Charlie/300) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "A"
$ws.Range("B1").Value = "B"
$ws.Range("C1").Value = "C"

# Leave the selection where the author apparently left it when saving.
$ws.Range("D6").Select()
